# Daily "cryptos" data refresh (GitHub Actions bot edit).
#
# All Coin/Link/Price/Volume(1h) cells in this sheet are stored as *text*
# (several Price values like "25.939.02" use dots as thousands separators
# and are not valid numbers at all; the rest are kept as text too so every
# cell in the column is formatted consistently and percentages keep their
# padding spaces). When writing through COM, a plain numeric-looking string
# (e.g. "4.23") gets auto-coerced to a real number by the Value/Value2
# setter, which would change the cell's type and drop formatting such as
# the shown decimals. Prefixing with a leading apostrophe forces Excel to
# store the entry as literal text (like a user typing '4.23 into the
# cell); resetting .Style back to "Normal" afterwards clears the transient
# "quote prefix" cell style that apostrophe-entry leaves behind, so styling
# stays identical to the untouched cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value2 = "'" + $text
    $cell.Style = "Normal"
}

# Bitcoin
Set-TextValue $ws.Range("D2") "25.934.10"
Set-TextValue $ws.Range("E2") "  +0.78%  "

# Ethereum
Set-TextValue $ws.Range("D3") "1.633.92"
Set-TextValue $ws.Range("E3") "  +0.38%  "

# TetherUSD
Set-TextValue $ws.Range("E4") "  +0.42%  "

# BNB
Set-TextValue $ws.Range("D5") "214.70"
Set-TextValue $ws.Range("E5") "  +0.19%  "

# XRP
Set-TextValue $ws.Range("E6") "  +0.98%  "

# USDC
Set-TextValue $ws.Range("E7") "  +0.35%  "

# Dogecoin
Set-TextValue $ws.Range("D9") "0.0632"
Set-TextValue $ws.Range("E9") "  +0.11%  "

# Solana
Set-TextValue $ws.Range("D10") "19.67"
Set-TextValue $ws.Range("E10") "  +1.02%  "

# TRON
Set-TextValue $ws.Range("E11") "  +0.12%  "

# WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.861.19"
Set-TextValue $ws.Range("E12") "  +0.39%  "

# Row 13 becomes WrappedEther (was Polkadot)
Set-TextValue $ws.Range("B13") "WrappedEther"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D13") "1.656.11"
Set-TextValue $ws.Range("E13") "  +1.80%  "

# Row 14 becomes Polkadot (was WrappedEther)
Set-TextValue $ws.Range("B14") "Polkadot"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "4.23"
Set-TextValue $ws.Range("E14") "  -0.45%  "

# Polygon
Set-TextValue $ws.Range("E15") "  -1.39%  "

# Litecoin
Set-TextValue $ws.Range("E16") "  +0.17%  "

# WrappedBTC
Set-TextValue $ws.Range("D18") "25.929.44"
Set-TextValue $ws.Range("E18") "  +0.72%  "

# Dai
Set-TextValue $ws.Range("E19") "  +0.41%  "

# BitcoinCash
Set-TextValue $ws.Range("D20") "193.15"
Set-TextValue $ws.Range("E20") "  +0.97%  "

# Uniswap
Set-TextValue $ws.Range("E21") "  -1.25%  "

# Chainlink
Set-TextValue $ws.Range("E23") "  -0.09%  "

# Toncoin
Set-TextValue $ws.Range("E24") "  -0.71%  "

# Monero
Set-TextValue $ws.Range("D25") "142.98"
Set-TextValue $ws.Range("E25") "  +0.39%  "

# BinanceUSD
Set-TextValue $ws.Range("E26") "  +0.21%  "

# Stellar
Set-TextValue $ws.Range("E27") "  +2.17%  "

# Cosmos
Set-TextValue $ws.Range("E28") "  +0.47%  "

# EthereumClassic
Set-TextValue $ws.Range("D29") "15.48"
Set-TextValue $ws.Range("E29") "  +0.04%  "

# PancakeSwap
Set-TextValue $ws.Range("E30") "  +0.31%  "

# Hedera
Set-TextValue $ws.Range("E31") "  +0.99%  "

# InternetComputer(DFINITY)
Set-TextValue $ws.Range("D32") "3.31"
Set-TextValue $ws.Range("E32") "  -0.43%  "

# Filecoin
Set-TextValue $ws.Range("E33") "  -0.10%  "

# LidoDAOToken
Set-TextValue $ws.Range("E34") "  -0.17%  "

# HuobiToken
Set-TextValue $ws.Range("E35") "  +2.15%  "

# ARBITRUM
Set-TextValue $ws.Range("E36") "  -0.45%  "

# Maker
Set-TextValue $ws.Range("D37") "1.137.69"
Set-TextValue $ws.Range("E37") "  -0.07%  "

# ImmutableX
Set-TextValue $ws.Range("E38") "  +1.72%  "

# MXToken
Set-TextValue $ws.Range("E39") "  -1.04%  "

# VeChain
Set-TextValue $ws.Range("E40") "  +0.60%  "

# PaxDollar
Set-TextValue $ws.Range("E41") "  +0.43%  "

# TrustWalletToken
Set-TextValue $ws.Range("E42") "  +0.15%  "

# FraxShare
Set-TextValue $ws.Range("E43") "  -1.45%  "

# Quant
Set-TextValue $ws.Range("D44") "99.20"
Set-TextValue $ws.Range("E44") "  -1.40%  "

# RocketPoolETH
Set-TextValue $ws.Range("D45") "1.770.28"
Set-TextValue $ws.Range("E45") "  +0.39%  "

# Row 46 becomes Aave (was BabyDogeCoin)
Set-TextValue $ws.Range("B46") "Aave"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "56.34"
Set-TextValue $ws.Range("E46") "  +2.28%  "

# Row 47 becomes Cronos (was Aave)
Set-TextValue $ws.Range("B47") "Cronos"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.0524"
Set-TextValue $ws.Range("E47") "  +2.67%  "

# Row 48 becomes RenderToken (was Cronos)
Set-TextValue $ws.Range("B48") "RenderToken"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D48") "1.46"
Set-TextValue $ws.Range("E48") "  +1.46%  "

# Row 49 becomes Mantle (was RenderToken)
Set-TextValue $ws.Range("B49") "Mantle"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.415"
Set-TextValue $ws.Range("E49") "  -0.39%  "

# Row 50 becomes EnergySwap (was Mantle)
Set-TextValue $ws.Range("B50") "EnergySwap"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.61"
Set-TextValue $ws.Range("E50") "  +1.45%  "

# Row 51 becomes Algorand (was EnergySwap)
Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0961"
Set-TextValue $ws.Range("E51") "  +0.80%  "
